$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original Text storage type (they were inline
# strings like "287.68" / "-1.16%") rather than being auto-parsed into
# numbers/percentages when we write the new values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "287.80"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.09%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.04"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.09%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.919"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.70%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07313"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.33%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.257"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "22.18%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.735"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.48%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.721"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.43%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9040"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.78%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09112"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "17.93%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1687"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.71%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08190"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.28%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03124"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.94%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09927"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.82%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001493"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.03%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005729"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.71%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.511"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.32%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.58%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3331"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.43%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.08%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.210"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "3.90%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-9.54%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04512"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.11%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001210"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.53%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-10.42%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "3.92%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01571"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-0.46%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04440"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.05%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007319"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.83%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009515"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-5.22%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1325"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.68%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002221"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "7.69%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.69%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006102"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.37%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.07%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.294"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "2.09%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.07%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.07%"
